# Actualizacion lista de precios lun 22/09/2025  9:40:33,37
#
# Adds four new "Ofertas" rows (with two brand-new article names that get
# appended to the shared-string table) and makes the "Ofertas" sheet the
# active tab, leaving cell D4 selected there.

$wb = $excel.ActiveWorkbook

$wsOfertas = $wb.Worksheets.Item("Ofertas")

# Fill in the new offers table on the "Ofertas" sheet.
$wsOfertas.Range("A1").Value = "Empresa 1"
$wsOfertas.Range("B1").Value = "Articulo 1"
$wsOfertas.Range("C1").Value = 2200

$wsOfertas.Range("A2").Value = "Empresa 1"
$wsOfertas.Range("B2").Value = "Articulo 3"
$wsOfertas.Range("C2").Value = 5490

$wsOfertas.Range("A3").Value = "Empresa 2"
$wsOfertas.Range("B3").Value = "Articulo 78"
$wsOfertas.Range("C3").Value = 1000

$wsOfertas.Range("A4").Value = "Empresa 1"
$wsOfertas.Range("B4").Value = "Articulo 9"
$wsOfertas.Range("C4").Value = 1200

# Make "Ofertas" the active/selected sheet and select D4 on it, which moves
# tabSelected from "Lista" to "Ofertas" and sets workbookView's activeTab.
$wsOfertas.Activate()
[void]$wsOfertas.Range("D4").Select()
